$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update region names and fan counts with the latest values (instead of mean values)
$ws.Range("A2").Value = "Attica (region), Greece"
$ws.Range("B2").Value = 115

$ws.Range("A3").Value = "Central Macedonia, Greece"
$ws.Range("B3").Value = 107

$ws.Range("A4").Value = "Thessaly, Greece"
$ws.Range("B4").Value = 26

$ws.Range("A5").Value = "Western Greece, Greece"
$ws.Range("B5").Value = 21

$ws.Range("A6").Value = "Eastern Macedonia and Thrace, Greece"
$ws.Range("B6").Value = 19

$ws.Range("A7").Value = "Central Greece (region), Greece"
$ws.Range("B7").Value = 18

$ws.Range("A8").Value = "Crete, Greece"
$ws.Range("B8").Value = 14

# Move the selection to A8, matching the saved workbook state
$ws.Range("A8").Select()
